$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace "kapil dev" with "Yuzvendra Chahal" in the points table (A12)
$ws.Range("A12").Value = "Yuzvendra Chahal"

# Column A width adjustment (as seen in new sheet XML: <col min="1" max="1" width="28" customWidth="1"/>)
$ws.Columns.Item(1).ColumnWidth = 27.14

# Update the active selection to K6, matching the updated sheetView
$ws.Range("K6").Select()
